# Append a new data row (row 25) to Sheet1, mirroring the formatting of the
# last existing row (row 24), as described by the diff / commit message
# ("Data set updated at Jan 3rd").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 25
$srcRow = 24

# Copy the formatting (number formats, styles) from the previous row down
# into the new row so the new data lines up visually with the rest of the
# table (same as row 24: A=date, B=day text, C/D numeric, E/F datetime,
# G/H/I numeric).
$ws.Range("A$srcRow`:I$srcRow").Copy() | Out-Null
$ws.Range("A$newRow`:I$newRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Now fill in the actual values for the new row.
$ws.Cells.Item($newRow, 1).Value = 43468.0            # A25 - Date (Jan 3, 2019)
$ws.Cells.Item($newRow, 2).Value = "Thu"               # B25 - Day
$ws.Cells.Item($newRow, 3).Value = 2.4                 # C25 - Stake
$ws.Cells.Item($newRow, 4).Value = 2.0                 # D25 - Toughness
$ws.Cells.Item($newRow, 5).Value = 43103.94305555556   # E25 - sesStart
$ws.Cells.Item($newRow, 6).Value = 43104.26944444444   # F25 - sesEnd
$ws.Cells.Item($newRow, 7).Value = 1200.0              # G25 - buyIn
$ws.Cells.Item($newRow, 8).Value = 0.0                 # H25 - promoBonus
$ws.Cells.Item($newRow, 9).Value = 1400.0              # I25 - cashOut
